# [Kadastro App] Yeni kayit eklendi: 3018
#
# Appends a new record (row 77) with Kayit No 3018 to both the "Kayitlar"
# master sheet and the "Erdemli" birim sheet (the workbook keeps a filtered
# copy of each "Birim" on its own tab, kept in sync with the master list).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

$newRow = @{
    A = "3018"
    B = "2025-09-11"
    C = "Erdemli"
    D = "1"
    E = "3B"
    F = "EMİNE ALANLI KIRCILI (K.Mühendisi), CEMAL TİMUROĞLU (K.Teknisyeni)"
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 77

    # Columns A, B, D hold digit-only / date-like text ("3018", "2025-09-11",
    # "1") that must stay stored as TEXT (matching every other row in the
    # sheet, which relies on the numberStoredAsText ignored-error), not be
    # auto-coerced into numbers/dates. Forcing the Text number format for the
    # assignment, then clearing formats again, keeps the cell's value as a
    # string while leaving its style back at the sheet default.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $newRow.A
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $newRow.B
    $ws.Cells.Item($r, 2).ClearFormats()

    $ws.Cells.Item($r, 3).Value = $newRow.C

    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $newRow.D
    $ws.Cells.Item($r, 4).ClearFormats()

    $ws.Cells.Item($r, 5).Value = $newRow.E
    $ws.Cells.Item($r, 6).Value = $newRow.F
}
